# Updates the "cryptos" sheet with refreshed Price/Volume(1h) figures and
# re-ranks a handful of coins (rows 41-46) to match the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.049.79"
$ws.Range("E2").Value = "  +7.69%  "
$ws.Range("D3").Value = "2.627.56"
$ws.Range("E3").Value = "  +8.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "185.60"
$ws.Range("E5").Value = "  +15.00%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "588.40"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.540"
$ws.Range("E8").Value = "  +4.99%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.203"
$ws.Range("E9").Value = "  +20.74%  "
$ws.Range("D10").Value = "2.622.42"
$ws.Range("E10").Value = "  +8.09%  "
$ws.Range("E11").Value = "  +0.34%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.365"
$ws.Range("E12").Value = "  +10.46%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.81"
$ws.Range("E13").Value = "  +4.64%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000191"
$ws.Range("E14").Value = "  +8.90%  "
$ws.Range("D15").Value = "74.124.09"
$ws.Range("E15").Value = "  +7.99%  "
$ws.Range("D16").Value = "3.112.75"
$ws.Range("E16").Value = "  +8.48%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.43"
$ws.Range("E17").Value = "  +14.38%  "
$ws.Range("D18").Value = "2.633.96"
$ws.Range("E18").Value = "  +8.54%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.07"
$ws.Range("E19").Value = "  +31.12%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.92"
$ws.Range("E20").Value = "  +13.82%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "376.47"
$ws.Range("E21").Value = "  +11.50%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.27"
$ws.Range("E22").Value = "  +18.10%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.12"
$ws.Range("E23").Value = "  +7.87%  "
$ws.Range("E24").Value = "  +0.04%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "70.21"
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("E26").Value = "  +14.37%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.42"
$ws.Range("E27").Value = "  +15.11%  "
$ws.Range("D28").Value = "2.747.93"
$ws.Range("E28").Value = "  +7.56%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.996"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "0.0₃0957"
$ws.Range("E30").Value = "  +17.79%  "
$ws.Range("E31").Value = "  +22.06%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.00"
$ws.Range("E32").Value = "  +12.56%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "512.37"
$ws.Range("E33").Value = "  +20.23%  "
$ws.Range("E34").Value = "  +9.38%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +16.55%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "160.00"
$ws.Range("E37").Value = "  +0.27%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "19.35"
$ws.Range("E38").Value = "  +8.00%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.33"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.95"
$ws.Range("E41").Value = "  +14.46%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.69"
$ws.Range("E42").Value = "  +12.84%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0929"
$ws.Range("E43").Value = "  +29.82%  "
$ws.Range("E44").Value = "  +10.17%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "159.76"
$ws.Range("E45").Value = "  +22.52%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.41"
$ws.Range("E46").Value = "  +18.87%  "
$ws.Range("E47").Value = "  +10.54%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "38.86"
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("E49").Value = "  +9.35%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.532"
$ws.Range("E50").Value = "  +10.86%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "20.73"
$ws.Range("E51").Value = "  +23.24%  "
